# Daily attendance processing - 2025-10-07 17:17:30
# Reorders the comma-separated "Recorded By" values in column G for the
# affected rows: the first name/email in the list is moved to the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,3,4,5,6,11,12,13,29,30,32,33,38,39,40,56,57,58,59,60,65,66,67,84,85,86,89,90,93,110,111,112,115,116,119,136,137,138,141,142,145)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $current = [string]$cell.Value2

    $parts = $current -split ',\s*'
    if ($parts.Count -gt 1) {
        $reordered = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ', '
        $cell.Value = $reordered
    }
}
